# historic_wind_speed.xlsx - "Updated weather and capacity data"
#
# The workbook already has monthly DateTime/Wind-speed pairs in columns A:B
# running through row 205 (2018-12-01). This continues the series with the
# 12 months of 2019 (rows 206-217, repeating the same 12-value seasonal
# cycle already present in the sheet) and formats a new column C (date
# style, no values yet) for rows 206-217, matching the author's in-progress
# edit. It also drops the four stale "_xlchart.v1.*" hidden defined names
# that referred to a deleted chart (#REF!).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the leftover hidden chart-reference defined names.
# ---------------------------------------------------------------------
$wb.Names.Item("_xlchart.v1.0").Delete()
$wb.Names.Item("_xlchart.v1.1").Delete()
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# ---------------------------------------------------------------------
# 2. Continue the DateTime / Wind speed series for another year
#    (2019-01-01 .. 2019-12-01), matching the repeating 12-month cycle
#    already present in column B.
# ---------------------------------------------------------------------
$newDates = @(43466, 43497, 43525, 43556, 43586, 43617, 43647, 43678, 43709, 43739, 43770, 43800)
$newSpeeds = @(998.59999999999991, 899.40000000000009, 893.39999999999986, 656.5, 596.1, 580.9, 568.99999999999989, 572.79999999999995, 586.20000000000005, 684.5, 743.6, 923.90000000000009)

$lastRow = 205
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $lastRow + 1 + $i

    # Carry the A-column date style (s="1") down onto the new row, then
    # overwrite with the actual date serial / wind-speed values.
    $ws.Range("A$lastRow").Copy($ws.Range("A$row"))
    $ws.Range("A$row").Value = $newDates[$i]
    $ws.Range("B$row").Value = $newSpeeds[$i]

    # Column C is formatted (same date style) but left empty for these
    # rows, same as the source edit.
    $ws.Range("A$lastRow").Copy($ws.Range("C$row"))
    $ws.Range("C$row").ClearContents()
}

# ---------------------------------------------------------------------
# 3. Leave the selection where the author left it: column C of the
#    newly-formatted rows.
# ---------------------------------------------------------------------
$ws.Range("C206:C217").Select()
